$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of planned-work data (rows 30 and 31)
$ws.Range("A30").Value = "Pokemon Sistemi DataBase'e entegre edilecek"
$ws.Range("B30").Value = "Pokemon"
$ws.Range("C30").Value = "Pokemon oluşturmak ScriptibleObject yerine JSON temelli olacak"
$ws.Range("D30").Value = "En Sona Doğru"

$ws.Range("A31").Value = "Shiny Pokemon"
$ws.Range("B31").Value = "Pokemon"
$ws.Range("C31").Value = "Shiny sistemi eklenecek"
$ws.Range("D31").Value = "Pokemon sistemi oturduktan sonra"

# Update the view state: scroll the window so row 22 is at the top, and
# leave the final selection on C36 (matches the saved workbook view).
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C36").Select()
